$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 5 = "TextBox 8": the "Release by Equalization (7 - 7.48 maf per
# year)" callout next to the Lower/Mid tier release label.
# ---------------------------------------------------------------------
$shape = $s.Shapes.Item(5)
$tr = $shape.TextFrame.TextRange

# Paragraph 1: "Release by Equalization" -> "Release"
$para1 = $tr.Paragraphs(1)
$para1.Text = "Release"

# Paragraph 2: "(7 - 7.48 maf per year)" -> "7 - 7.48 maf per year"
# (drop the leading "(" and trailing ")", keep the "maf" run intact so
# its spell-check flag/formatting survives unchanged)
$para2 = $tr.Paragraphs(2)
$para2.Characters(1, 1).Text = ""
$para2 = $tr.Paragraphs(2)
$para2.Characters($para2.Text.Length, 1).Text = ""

# New paragraph 3: "(Lower - Mid tier)"
$tr.InsertAfter([char]13 + "(Lower ") | Out-Null
$tr.InsertAfter([char]8211 + " ") | Out-Null
$tr.InsertAfter("Mid tier)") | Out-Null

# Resize/reposition the textbox to match its new (autofit) content box;
# height grows to fit the extra line while the box stays horizontally
# centered on its original midpoint.
$shape.Left = 86.0716248031496
$shape.Width = 188.7403649606299
$shape.Height = 79.97351456692914

# ---------------------------------------------------------------------
# Shape 11 = "Straight Arrow Connector 24": glued to the shape above
# (endCxn id=9/idx=3), so its bounding box has to follow the textbox's
# new position/size.
# ---------------------------------------------------------------------
$connector = $s.Shapes.Item(11)
$connector.Left = 274.8119397637795
$connector.Top = 268.06997125984253
$connector.Width = 89.28973503937009
$connector.Height = 63.15690039370079
